$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'28.515.95"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  +1.50%  '

$ws.Cells.Item(3, 4).Value = "'1.823.80"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +1.32%  '

$ws.Cells.Item(4, 4).Value = "'1.005"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  +0.37%  '

$ws.Cells.Item(5, 4).Value = "'315.82"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.32%  '

$ws.Cells.Item(6, 4).Value = "'1.004"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +0.23%  '

$ws.Cells.Item(7, 4).Value = "'0.5150"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -5.61%  '

$ws.Cells.Item(8, 4).Value = "'0.3919"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  +3.55%  '

$ws.Cells.Item(9, 4).Value = "'0.07695"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +3.10%  '

$ws.Cells.Item(10, 4).Value = "'42.04"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +0.20%  '

$ws.Cells.Item(11, 4).Value = "'1.110"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +1.58%  '

$ws.Cells.Item(12, 4).Value = "'20.97"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  +2.72%  '

$ws.Cells.Item(13, 2).Value = 'BinanceUSD'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(13, 4).Value = "'1.006"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +0.41%  '

$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).Value = "'6.272"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +1.34%  '

$ws.Cells.Item(15, 4).Value = "'7.538"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +2.64%  '

$ws.Cells.Item(16, 4).Value = "'1.824.62"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +1.42%  '

$ws.Cells.Item(17, 4).Value = "'92.93"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +3.84%  '

$ws.Cells.Item(18, 4).Value = "'0.00001077"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +1.28%  '

$ws.Cells.Item(19, 4).Value = "'0.06597"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +0.79%  '

$ws.Cells.Item(20, 4).Value = "'17.66"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +1.52%  '

$ws.Cells.Item(21, 4).Value = "'1.003"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +0.24%  '

$ws.Cells.Item(22, 4).Value = "'6.061"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +2.17%  '

$ws.Cells.Item(23, 4).Value = "'28.533.70"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +1.43%  '

$ws.Cells.Item(24, 4).Value = "'11.11"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.74%  '

$ws.Cells.Item(25, 4).Value = "'2.242"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +7.48%  '

$ws.Cells.Item(26, 4).Value = "'20.62"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +1.02%  '

$ws.Cells.Item(27, 4).Value = "'156.29"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +0.44%  '

$ws.Cells.Item(28, 4).Value = "'2.037.30"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.55%  '

$ws.Cells.Item(29, 4).Value = "'2.404"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +3.59%  '

$ws.Cells.Item(30, 4).Value = "'124.82"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +2.56%  '

$ws.Cells.Item(31, 4).Value = "'1.132"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +1.51%  '

$ws.Cells.Item(32, 4).Value = "'0.1110"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -0.50%  '

$ws.Cells.Item(33, 4).Value = "'5.666"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +2.02%  '

$ws.Cells.Item(34, 4).Value = "'3.661"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -0.47%  '

$ws.Cells.Item(35, 4).Value = "'0.07220"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +4.79%  '

$ws.Cells.Item(36, 4).Value = "'0.2240"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +0.85%  '

$ws.Cells.Item(37, 4).Value = "'0.02335"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +1.94%  '

$ws.Cells.Item(38, 4).Value = "'8.876"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +5.04%  '

$ws.Cells.Item(39, 4).Value = "'5.140"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +1.03%  '

$ws.Cells.Item(40, 4).Value = "'11.29"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +0.87%  '

$ws.Cells.Item(41, 4).Value = "'0.6234"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +1.26%  '

$ws.Cells.Item(42, 4).Value = "'1.186"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +1.15%  '

$ws.Cells.Item(43, 4).Value = "'1.003"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.27%  '

$ws.Cells.Item(44, 5).Value = '  -1.60%  '

$ws.Cells.Item(45, 4).Value = "'13.43"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +1.02%  '

$ws.Cells.Item(46, 2).Value = 'PancakeSwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(46, 4).Value = "'3.717"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +0.96%  '

$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).Value = "'0.5892"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +2.69%  '

$ws.Cells.Item(48, 4).Value = "'124.89"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +0.20%  '

$ws.Cells.Item(49, 4).Value = "'1.979"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +3.18%  '

$ws.Cells.Item(50, 4).Value = "'1.186"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +0.35%  '

$ws.Cells.Item(51, 4).Value = "'0.06940"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +1.87%  '
